# "Substituindo planilha por versão limpa e corrigida"
# Rename the only sheet, move the active-cell selection, and tighten up
# the data-column widths (B, F:M) to match the cleaned-up layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename "Planilha1" -> "BASE"
$ws.Name = "BASE"

# 2. Column width adjustments.
#    Target XML `width` values (character units) from the cleaned sheet,
#    converted to the `ColumnWidth` COM property (which is offset by the
#    standard 5/6-character cell-margin before Excel serialises it back
#    out to the `width` attribute).
$ws.Columns.Item(2).ColumnWidth  = 15.053385416666666   # B: 16.77734375 -> 15.88671875
$ws.Columns.Item(6).ColumnWidth  = 10.830729166666666   # F: 13.6640625  -> 11.6640625
$ws.Columns.Item(7).ColumnWidth  = 11.053385416666666   # G: 11          -> 11.88671875
$ws.Columns.Item(8).ColumnWidth  = 11.053385416666666   # H: 11.44140625 -> 11.88671875
$ws.Columns.Item(9).ColumnWidth  = 11.830729166666666   # I: 11.33203125 -> 12.6640625
$ws.Columns.Item(10).ColumnWidth = 11.830729166666666   # J: 10.33203125 -> 12.6640625
$ws.Columns.Item(11).ColumnWidth = 11.830729166666666   # K: 11.5546875  -> 12.6640625
$ws.Columns.Item(12).ColumnWidth = 11.830729166666666   # L: 10.5546875  -> 12.6640625
$ws.Columns.Item(13).ColumnWidth = 10.498697916666666   # M: 12.44140625 -> 11.33203125

# 3. Move the active-cell selection from N4 to N6.
$ws.Range("N6").Select()
